# Apply the "Working test of DB and PR interfaces" edit to RoperSpreadSheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 4: quick single-row smoke test of the intake columns ---
$ws.Range("B4").Value = "test"
$ws.Range("C4").Value = "test"
$ws.Range("D4").Value = Get-Date -Year 1987 -Month 6 -Day 5 -Hour 0 -Minute 0 -Second 0
$ws.Range("F4").Value = "W"
$ws.Range("G4").Value = "M"

# --- Rows 33 & 34: a second (duplicated) test record ---
$dob2 = Get-Date -Year 1969 -Month 12 -Day 31 -Hour 19 -Minute 0 -Second 0
foreach ($r in 33, 34) {
    $ws.Range("B$r").Value = "Franks"
    $ws.Range("C$r").Value = "Kevin"
    $ws.Range("D$r").Value = $dob2
    $ws.Range("F$r").Value = "Aliens"
    $ws.Range("G$r").Value = "C"
}

# --- Column D (DOB) widened to fit the new dates ---
$ws.Columns.Item(4).ColumnWidth = 9.721354166666666

# --- Rename the hidden WorksheetConnection defined names (xlcn table refs) ---
$wb.Names.Item("_xlcn.WorksheetConnection_newdatabase.xlsxTable1").Name = "_xlcn.WorksheetConnection_newdatabase.xlsxTable11"
$wb.Names.Item("_xlcn.WorksheetConnection_newdatabase.xlsxTable2").Name = "_xlcn.WorksheetConnection_newdatabase.xlsxTable21"

# --- Selection / view state left pointing at the newly entered row ---
$ws.Range("D34").Select() | Out-Null
